# Simulator full-month coverage, persist logs, fix employees
$wb = $excel.ActiveWorkbook

$wsTime = $wb.Worksheets.Item("Weekly Timesheet")
$wsSchema = $wb.Worksheets.Item("Jason Schema")

# ---- Fix employee/client names (rows 2-6 correspond to 2026-01-19 .. 2026-01-23) ----
$clients = @("Tubergen", "Hewett", "Durfee", "Markfield", "Corr")

for ($i = 0; $i -lt $clients.Length; $i++) {
    $row = 2 + $i
    $wsTime.Range("B$row").Value = $clients[$i]
    $wsSchema.Range("D$row").Value = $clients[$i]
}

# ---- Simulator full-month coverage: hours 8 -> 9, rate 0 -> 150, total 0 -> 1350 ----
for ($row = 2; $row -le 6; $row++) {
    $wsTime.Range("C$row").Value = 9
    $wsTime.Range("E$row").Value = 150
    $wsTime.Range("F$row").Value = 1350

    $wsSchema.Range("E$row").Value = 9
    $wsSchema.Range("F$row").Value = 150
    $wsSchema.Range("G$row").Value = 1350
}

# ---- Weekly Timesheet subtotal / totals ----
$wsTime.Range("C8").Value = 45
$wsTime.Range("F8").Value = 6750
$wsTime.Range("F12").Value = 6750
$wsTime.Range("F13").Value = 6750

# ---- Reg/OT summary text (row 8, column D: "Reg: 40 / OT: 0" -> "Reg: 45 / OT: 0") ----
$wsTime.Range("D8").Value = "Reg: 45 / OT: 0"

# ---- Fix employee ID (appears once per data row on Jason Schema) ----
for ($row = 2; $row -le 6; $row++) {
    $wsSchema.Range("B$row").Value = "emp_35u1tnme"
}
